$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-02-19 Wednesday" "2025-02-20 Thursday"

Replace-Text "652÷7=93, 1" "445÷2=222, 1"
Replace-Text "949÷4=237, 1" "779÷7=111, 2"
Replace-Text "472÷9=52, 4" "132÷9=14, 6"
Replace-Text "695÷7=99, 2" "798÷7=114, 0"
Replace-Text "999÷5=199, 4" "640÷3=213, 1"

Replace-Text "514÷2=257, 0" "297÷4=74, 1"
Replace-Text "630÷5=126, 0" "523÷7=74, 5"
Replace-Text "481÷7=68, 5" "363÷3=121, 0"
Replace-Text "464÷3=154, 2" "303÷7=43, 2"
Replace-Text "117÷7=16, 5" "591÷4=147, 3"

Replace-Text "769÷5=153, 4" "457÷6=76, 1"
Replace-Text "711÷2=355, 1" "116÷2=58, 0"
Replace-Text "591÷8=73, 7" "696÷5=139, 1"
Replace-Text "895÷5=179, 0" "432÷7=61, 5"
Replace-Text "884÷9=98, 2" "883÷4=220, 3"

Replace-Text "838÷5=167, 3" "430÷2=215, 0"
Replace-Text "940÷2=470, 0" "870÷7=124, 2"
Replace-Text "900÷7=128, 4" "944÷6=157, 2"
Replace-Text "696÷7=99, 3" "449÷9=49, 8"
Replace-Text "971÷8=121, 3" "439÷9=48, 7"

Replace-Text "284÷7=40, 4" "252÷5=50, 2"
Replace-Text "965÷3=321, 2" "165÷2=82, 1"
Replace-Text "679÷5=135, 4" "415÷2=207, 1"
Replace-Text "900÷2=450, 0" "305÷3=101, 2"
Replace-Text "901÷8=112, 5" "759÷7=108, 3"
